$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.966.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.811.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.804.80"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.441.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.806.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.112.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.00%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.56%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.62"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.61"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.42"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.35%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.02"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.55"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.866.66"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0353"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.20"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.47"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.27%  "
